$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking strings
# (e.g. "209.83", "1.78") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "87.858.52"
$ws.Range("E2").Value = "  -2.85%  "
$ws.Range("D3").Value = "3.053.35"
$ws.Range("E3").Value = "  -4.89%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "209.83"
$ws.Range("E5").Value = "  -3.98%  "
$ws.Range("D6").Value = "615.73"
$ws.Range("E6").Value = "  -5.63%  "
$ws.Range("D7").Value = "0.368"
$ws.Range("E7").Value = "  -7.64%  "
$ws.Range("D8").Value = "0.794"
$ws.Range("E8").Value = "  +13.99%  "
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("D10").Value = "3.051.77"
$ws.Range("E10").Value = "  -4.92%  "
$ws.Range("D11").Value = "0.592"
$ws.Range("E11").Value = "  +2.15%  "
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("E13").Value = "  -7.95%  "
$ws.Range("E14").Value = "  -2.97%  "
$ws.Range("D15").Value = "87.770.39"
$ws.Range("E15").Value = "  -2.52%  "
$ws.Range("D16").Value = "3.622.21"
$ws.Range("E16").Value = "  -4.59%  "
$ws.Range("D17").Value = "31.72"
$ws.Range("E17").Value = "  -5.64%  "
$ws.Range("D18").Value = "3.061.97"
$ws.Range("E18").Value = "  -4.56%  "
$ws.Range("D19").Value = "3.27"
$ws.Range("E19").Value = "  -4.65%  "
$ws.Range("D20").Value = "0.0000198"
$ws.Range("E20").Value = "  -12.88%  "
$ws.Range("D21").Value = "13.20"
$ws.Range("E21").Value = "  -3.07%  "
$ws.Range("D22").Value = "419.21"
$ws.Range("E22").Value = "  -4.99%  "
$ws.Range("D23").Value = "8.10"
$ws.Range("E23").Value = "  -6.90%  "
$ws.Range("D24").Value = "4.88"
$ws.Range("E24").Value = "  -4.73%  "
$ws.Range("D25").Value = "5.44"
$ws.Range("E25").Value = "  +1.76%  "
$ws.Range("D26").Value = "11.73"
$ws.Range("E26").Value = "  -2.65%  "
$ws.Range("D27").Value = "81.80"
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("D28").Value = "3.231.99"
$ws.Range("E28").Value = "  -4.05%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Value = "1.09"
$ws.Range("E30").Value = "  +9.31%  "
$ws.Range("D31").Value = "0.171"
$ws.Range("E31").Value = "  +6.38%  "
$ws.Range("D32").Value = "8.01"
$ws.Range("E32").Value = "  -6.36%  "
$ws.Range("D33").Value = "502.70"
$ws.Range("E33").Value = "  -8.71%  "
$ws.Range("E34").Value = "  -12.24%  "
$ws.Range("D35").Value = "6.69"
$ws.Range("E35").Value = "  -6.25%  "
$ws.Range("B36").Value = "PancakeSwap"
$ws.Range("C36").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D36").Value = "1.78"
$ws.Range("E36").Value = "  -8.47%  "
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").Value = "1.23"
$ws.Range("E37").Value = "  -6.48%  "
$ws.Range("D38").Value = "22.12"
$ws.Range("E38").Value = "  -2.30%  "
$ws.Range("E39").Value = "  -1.02%  "
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("E43").Value = "  -4.80%  "
$ws.Range("D44").Value = "147.75"
$ws.Range("E44").Value = "  +0.74%  "
$ws.Range("E45").Value = "  -8.49%  "
$ws.Range("E46").Value = "  +4.98%  "
$ws.Range("E47").Value = "  -4.03%  "
$ws.Range("D48").Value = "0.0683"
$ws.Range("E48").Value = "  +12.28%  "
$ws.Range("D49").Value = "156.73"
$ws.Range("E49").Value = "  -10.30%  "
$ws.Range("E50").Value = "  -5.80%  "
$ws.Range("E51").Value = "  -9.04%  "

# Restore default style on column D so no stray number-format
# style is left behind on cells (matches original formatting).
$ws.Range("D2:D51").Style = "Normal"
